$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 219, shifting existing rows 219:233 down to 220:234
$ws.Rows.Item(219).Insert()

# Populate the new row 219 with the new record (dated 2022-06-02 = serial 44714)
$ws.Range("A219").Value = 10
$ws.Range("B219").Value = "Vega Modelo de Temuco"
$ws.Range("C219").Value = "La Araucanía"
$ws.Range("D219").Value = 44714
$ws.Range("E219").Value = 9
$ws.Range("F219").Value = 100112043
$ws.Range("G219").Value = "Pepino dulce"
$ws.Range("H219").Value = "Cultivar IV Región"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 65
$ws.Range("K219").Value = 17000
$ws.Range("L219").Value = 17000
$ws.Range("M219").Value = 17000
$ws.Range("N219").Value = "$/bandeja 18 kilos"
$ws.Range("O219").Value = "Provincia de Limarí"
$ws.Range("P219").Value = 944
$ws.Range("Q219").Value = 18
$ws.Range("R219").Value = "Hortaliza"
